$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3 ("Bahasa Daerah"), shifting rows 3:26 down
# to 4:27 and making room for a new "BK" (Bimbingan Konseling) subject row.
$ws.Rows.Item(3).Insert()

# Seed the new row 3 from the row that landed just below it (a full copy
# keeps every column, including the blank "Deskripsi" cell, present) and
# then overwrite the columns that actually hold new data.
$ws.Range("A4:E4").Copy($ws.Range("A3:E3"))

$ws.Range("A3").Value = "BK"
$ws.Range("B3").Value = "Bimbingan Konseling"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "SEMUA"
